$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.Value2 = $text
    $rng.Style = $origStyle
}

Set-CellText "D2" '67.234.48'
Set-CellText "E2" '  +3.29%  '

Set-CellText "D3" '3.467.26'
Set-CellText "E3" '  +2.75%  '

Set-CellText "D4" '1.00'
Set-CellText "E4" '  +0.09%  '

Set-CellText "D5" '584.56'
Set-CellText "E5" '  +5.06%  '

Set-CellText "D6" '191.22'
Set-CellText "E6" '  +9.34%  '

Set-CellText "D7" '0.634'
Set-CellText "E7" '  +0.50%  '

Set-CellText "D8" '3.460.19'
Set-CellText "E8" '  +2.89%  '

Set-CellText "D9" '1.00'
Set-CellText "E9" '  +0.01%  '

Set-CellText "D10" '0.173'
Set-CellText "E10" '  +0.09%  '

Set-CellText "D11" '0.650'
Set-CellText "E11" '  +2.03%  '

Set-CellText "D12" '57.78'
Set-CellText "E12" '  +7.57%  '

Set-CellText "E13" '  -0.01%  '

Set-CellText "D14" '9.53'
Set-CellText "E14" '  +3.47%  '

Set-CellText "D15" '4.009.66'
Set-CellText "E15" '  +2.62%  '

Set-CellText "D16" '18.95'
Set-CellText "E16" '  +3.64%  '

Set-CellText "D17" '3.470.71'
Set-CellText "E17" '  +2.91%  '

Set-CellText "D18" '67.181.93'
Set-CellText "E18" '  +3.61%  '

Set-CellText "E19" '  +2.30%  '

Set-CellText "E20" '  +1.07%  '

Set-CellText "E21" '  +3.07%  '

Set-CellText "D22" '484.49'
Set-CellText "E22" '  +6.56%  '

Set-CellText "B23" 'Toncoin'
Set-CellText "C23" 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-CellText "D23" '5.35'
Set-CellText "E23" '  +8.77%  '

Set-CellText "B24" 'InternetComputer(DFINITY)'
Set-CellText "C24" 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-CellText "D24" '16.78'
Set-CellText "E24" '  +19.21%  '

Set-CellText "E25" '  +7.47%  '

Set-CellText "D26" '90.47'
Set-CellText "E26" '  +3.37%  '

Set-CellText "D27" '3.01'
Set-CellText "E27" '  +4.58%  '

Set-CellText "E28" '  +2.83%  '

Set-CellText "D29" '9.10'
Set-CellText "E29" '  +4.32%  '

Set-CellText "D30" '31.45'
Set-CellText "E30" '  +0.79%  '

Set-CellText "D31" '7.44'
Set-CellText "E31" '  +13.71%  '

Set-CellText "B32" 'Bittensor'
Set-CellText "C32" 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-CellText "D32" '603.50'
Set-CellText "E32" '  +4.43%  '

Set-CellText "B33" 'Cosmos'
Set-CellText "C33" 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-CellText "D33" '11.86'
Set-CellText "E33" '  +3.50%  '

Set-CellText "D34" '64.45'
Set-CellText "E34" '  +2.18%  '

Set-CellText "D35" '0.113'
Set-CellText "E35" '  +4.66%  '

Set-CellText "E36" '  +5.99%  '

Set-CellText "E37" '  -0.09%  '

Set-CellText "D38" '37.52'
Set-CellText "E38" '  +5.08%  '

Set-CellText "D39" '0.391'
Set-CellText "E39" '  +5.05%  '

Set-CellText "D40" '3.50'
Set-CellText "E40" '  -4.70%  '

Set-CellText "D41" '0.0₃0760'
Set-CellText "E41" '  +2.55%  '

Set-CellText "D42" '3.232.69'
Set-CellText "E42" '  +4.55%  '

Set-CellText "D43" '2.95'
Set-CellText "E43" '  +6.48%  '

Set-CellText "B44" 'dogwifhat'
Set-CellText "C44" 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-CellText "D44" '2.95'
Set-CellText "E44" '  +31.04%  '

Set-CellText "B45" 'VeChain'
Set-CellText "C45" 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-CellText "D45" '0.0434'
Set-CellText "E45" '  +4.20%  '

Set-CellText "D46" '2.57'
Set-CellText "E46" '  +4.38%  '

Set-CellText "E47" '  +1.49%  '

Set-CellText "E48" '  +1.18%  '

Set-CellText "B49" 'THORChain'
Set-CellText "C49" 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-CellText "D49" '8.80'
Set-CellText "E49" '  +6.27%  '

Set-CellText "B50" 'FirstDigitalUSD'
Set-CellText "C50" 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-CellText "D50" '1.00'
Set-CellText "E50" '  +0.15%  '

Set-CellText "B51" 'LidoDAOToken'
Set-CellText "C51" 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-CellText "D51" '3.23'
Set-CellText "E51" '  +7.95%  '
